$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(38, 1).Value = 'P76'
$ws.Cells.Item(38, 2).Value = 1754583810409
$ws.Cells.Item(38, 3).Value = '2025-08-07 18:23:30.409'
$ws.Cells.Item(38, 4).Value = 'Apply Button Pressed'
$ws.Cells.Item(38, 5).Value = '{''taskMode'': ''First: Yellow'', ''taskOrder'': ''Free''}'

$ws.Cells.Item(39, 1).Value = 'P76'
$ws.Cells.Item(39, 2).Value = 1754583810907
$ws.Cells.Item(39, 3).Value = '2025-08-07 18:23:30.907'
$ws.Cells.Item(39, 4).Value = 'Initialize Robot Button Pressed'
$ws.Cells.Item(39, 5).Value = '{}'

$ws.Cells.Item(40, 1).Value = 'P76'
$ws.Cells.Item(40, 2).Value = 1754583816966
$ws.Cells.Item(40, 3).Value = '2025-08-07 18:23:36.966'
$ws.Cells.Item(40, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(40, 5).Value = '{''taskId'': ''1'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(41, 1).Value = 'P76'
$ws.Cells.Item(41, 2).Value = 1754583817031
$ws.Cells.Item(41, 3).Value = '2025-08-07 18:23:37.031'
$ws.Cells.Item(41, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(41, 5).Value = '{''taskId'': ''2'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(42, 1).Value = 'P76'
$ws.Cells.Item(42, 2).Value = 1754583817100
$ws.Cells.Item(42, 3).Value = '2025-08-07 18:23:37.100'
$ws.Cells.Item(42, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(42, 5).Value = '{''taskId'': ''3'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(43, 1).Value = 'P76'
$ws.Cells.Item(43, 2).Value = 1754583817165
$ws.Cells.Item(43, 3).Value = '2025-08-07 18:23:37.165'
$ws.Cells.Item(43, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(43, 5).Value = '{''taskId'': ''5'', ''assignedTo'': ''Robot'', ''sliderValue'': 10}'

$ws.Cells.Item(44, 1).Value = 'P76'
$ws.Cells.Item(44, 2).Value = 1754583817307
$ws.Cells.Item(44, 3).Value = '2025-08-07 18:23:37.307'
$ws.Cells.Item(44, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(44, 5).Value = '{''taskId'': ''6'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(45, 1).Value = 'P76'
$ws.Cells.Item(45, 2).Value = 1754583817428
$ws.Cells.Item(45, 3).Value = '2025-08-07 18:23:37.428'
$ws.Cells.Item(45, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(45, 5).Value = '{''taskId'': ''8'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(46, 1).Value = 'P76'
$ws.Cells.Item(46, 2).Value = 1754583817542
$ws.Cells.Item(46, 3).Value = '2025-08-07 18:23:37.542'
$ws.Cells.Item(46, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(46, 5).Value = '{''taskId'': ''9'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(47, 1).Value = 'P76'
$ws.Cells.Item(47, 2).Value = 1754583817642
$ws.Cells.Item(47, 3).Value = '2025-08-07 18:23:37.642'
$ws.Cells.Item(47, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(47, 5).Value = '{''taskId'': ''10'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(48, 1).Value = 'P76'
$ws.Cells.Item(48, 2).Value = 1754583817942
$ws.Cells.Item(48, 3).Value = '2025-08-07 18:23:37.942'
$ws.Cells.Item(48, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(48, 5).Value = '{''taskId'': ''12'', ''assignedTo'': ''Human'', ''sliderValue'': 2}'

$ws.Cells.Item(49, 1).Value = 'P76'
$ws.Cells.Item(49, 2).Value = 1754583818818
$ws.Cells.Item(49, 3).Value = '2025-08-07 18:23:38.818'
$ws.Cells.Item(49, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(49, 5).Value = '{''taskId'': ''14'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(50, 1).Value = 'P76'
$ws.Cells.Item(50, 2).Value = 1754583819569
$ws.Cells.Item(50, 3).Value = '2025-08-07 18:23:39.569'
$ws.Cells.Item(50, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(50, 5).Value = '{''taskId'': ''15'', ''assignedTo'': ''Human'', ''sliderValue'': 1}'

$ws.Cells.Item(51, 1).Value = 'P76'
$ws.Cells.Item(51, 2).Value = 1754583820610
$ws.Cells.Item(51, 3).Value = '2025-08-07 18:23:40.610'
$ws.Cells.Item(51, 4).Value = 'Task Allocation Changed'
$ws.Cells.Item(51, 5).Value = '{''taskId'': ''18'', ''assignedTo'': ''Human'', ''sliderValue'': 0}'

$ws.Cells.Item(52, 1).Value = 'P76'
$ws.Cells.Item(52, 2).Value = 1754583822249
$ws.Cells.Item(52, 3).Value = '2025-08-07 18:23:42.249'
$ws.Cells.Item(52, 4).Value = 'Start Button Pressed'
$ws.Cells.Item(52, 5).Value = '{}'

$ws.Cells.Item(53, 1).Value = 'Robot'
$ws.Cells.Item(53, 2).Value = 1754583839907
$ws.Cells.Item(53, 3).Value = '2025-08-07 18:23:59.907'
$ws.Cells.Item(53, 4).Value = 'Robot Task Completed'
$ws.Cells.Item(53, 5).Value = '{''task_name'': ''Bridge_triangle_roof'', ''urp_name'': ''Bridge_triangle_roof''}'
